$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.373.22"
$ws.Range("D3").Value = "1.879.55"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'0.7214"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").Value = "'243.06"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D8").Value = "'0.08015"
$ws.Range("E8").Value = "  +2.66%  "
$ws.Range("D9").Value = "'0.3140"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "'24.96"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "'0.08155"
$ws.Range("E11").Value = "  -3.20%  "
$ws.Range("D12").Value = "1.883.61"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "'94.58"
$ws.Range("E13").Value = "  +3.83%  "
$ws.Range("D14").Value = "'5.229"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "'0.7116"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("E16").Value = "  +5.64%  "
$ws.Range("D17").Value = "'0.000008471"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value = "29.375.76"
$ws.Range("D19").Value = "'244.23"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").Value = "'13.30"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "2.127.35"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'7.744"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'0.1604"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").Value = "'162.69"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'9.036"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'18.50"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "'1.505"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'4.399"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'4.281"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").Value = "'1.231"
$ws.Range("E32").Value = "  -5.42%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "'1.938"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").Value = "'0.7641"
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("D36").Value = "'1.177"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'2.699"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").Value = "'0.01871"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("D39").Value = "1.261.32"
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").Value = "'2.765"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "'6.438"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").Value = "'113.30"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.9051"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'74.15"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("E45").Value = "  +5.85%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "2.025.86"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").Value = "'1.801"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "'0.5201"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'9.482"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "'0.4335"
$ws.Range("E51").Value = "  +0.04%  "
